$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update page count for "The Passionate Programmer" (row 13)
$ws.Range("C13").Value = 133

# Force recalculation so the SUM formula in E7 picks up the new total
$excel.Calculate()

# Move the active selection to D25 (matches the saved sheetView state)
$ws.Range("D25").Select()
